# Department API workbook edit
# - Rename sheets 2-5 from "학교 공지 ..." (school notice) to "학과 공지 ..." (department notice)
#   to match the workbook's actual subject (Department API).
# - Update the active selection on the last sheet ("학과 공지 삭제") from K6 to B24.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(2).Name = "학과 공지 상세 조회"
$wb.Worksheets.Item(3).Name = "학과 공지 작성"
$wb.Worksheets.Item(4).Name = "학과 공지 수정"
$wb.Worksheets.Item(5).Name = "학과 공지 삭제"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$ws5.Range("B24").Select()
